$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds one 4-trial block (rows 2-5): 1+2, 1+3, 3+4, 2+4.
# The task now needs 20 trials per block, i.e. that same 4-row cycle
# repeated 5 times (rows 2-21).
$condLabels = "1+2", "1+3", "3+4", "2+4"
$bVals = 1, 1, 3, 2
$cVals = 2, 3, 4, 4

$row = 6
for ($block = 1; $block -le 4; $block++) {
    for ($i = 0; $i -lt 4; $i++) {
        $ws.Cells.Item($row, 1).Value = $condLabels[$i]
        $ws.Cells.Item($row, 2).Value = $bVals[$i]
        $ws.Cells.Item($row, 3).Value = $cVals[$i]
        $row++
    }
}

# The "3+4" condition rows (A4, A8, A12, A16, A20) carry their own cell
# format (same font, just flagged as explicitly applied).
foreach ($r in 4, 8, 12, 16, 20) {
    $ws.Cells.Item($r, 1).Font.ThemeColor = 1
}

$ws.Range("D19").Select() | Out-Null
